$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 68. This pushes the current totals row
# (old 68) down to 69 and the footer row (old 69) down to 70, and shifts the
# related merged-cell ranges automatically.
$ws.Rows("68:68").Insert()

# Row 67 currently still holds the "محلول ملح" line (values carried over from
# before the insert). Duplicate it down into the newly created row 68 so that
# item keeps its place right after the new product line we are about to
# insert into row 67.
$ws.Range("A67:Q67").Copy($ws.Range("A68:Q68"))
$ws.Rows("68:68").RowHeight = 24.75

# Fix the sequence number + total-group cell in the duplicated row (row 68
# becomes item #62).
$ws.Cells.Item(68, 1).Value = 62

# Recreate the merged cells for the new row 68 (it had none before, since it
# used to be an empty spacer row).
$ws.Range("A68:B68").Merge()
$ws.Range("C68:G68").Merge()
$ws.Range("H68:K68").Merge()
$ws.Range("L68:M68").Merge()
$ws.Range("N68:O68").Merge()

# Now overwrite row 67 with the data for the newly stocked product
# "كريم فيرند لافلي الكبير". A, B, C, D-G stay put (same item cell / merges);
# only the quantity columns change.
$h67 = $ws.Cells.Item(67, 8)
$h67.NumberFormat = "@"
$h67.Value = "8:0"

$n67 = $ws.Cells.Item(67, 14)
$n67.NumberFormat = "@"
$n67.Value = "35.00"

$p67 = $ws.Cells.Item(67, 16)
$p67.NumberFormat = "@"
$p67.Value = "35.0000"
$p67.NumberFormat = "0.00"

$q67 = $ws.Cells.Item(67, 17)
$q67.NumberFormat = "@"
$q67.Value = "1:0"

# C67 keeps referencing the shared string that used to mean "محلول ملح" but
# now needs to say "كريم فيرند لافلي الكبير" instead (the old row's text was
# already duplicated into row 68 above, so it's safe to overwrite here).
$c67 = $ws.Cells.Item(67, 3)
$c67.NumberFormat = "@"
$c67.Value = "كريم فيرند لافلي الكبير"

# The daily total (now on row 69) grows by the new item's sale value, and
# the row grows slightly taller to match the committed workbook.
$ws.Cells.Item(69, 16).Value = 3824.5949999999998
$ws.Rows("69:69").RowHeight = 25.5
